# billing.xlsx edit: "mau upload and billing report"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date label updates (top banner + row 7) before the row shuffle ---
$ws.Range("A3").Value = "23 jun 2023"
$ws.Range("C7").Value = "600     Monthly Active Users @`$601/month"
$ws.Range("D7").Value = "23 jun 2023"

# --- Remove the old "800 Agent Seats" row and the old "Platform Support" row.
# Deleting row 9 twice removes both (the second delete removes what used to be
# row 10 once everything below row 9 has shifted up). This brings the old
# "WHATSAPP CONVERSATIONS" section header (previously row 11) up to row 9,
# and shifts the conversation-fee rows up to rows 10-14, matching the target
# layout/styles exactly (Excel's row delete keeps formatting of the rows that
# shift into place).
$ws.Rows(9).Delete()
$ws.Rows(9).Delete()

# --- Row 8 becomes the new "Platform Support" line item ---
$ws.Range("C8").Value = "Platform Support"
$ws.Range("D8").Value = "23 jun 2023"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = 100

# --- Row 9 is already "WHATSAPP CONVERSATIONS" after the shift; nothing to edit ---

# --- Row 10: "Fee Conversation/Month" with new qty/fee values ---
$ws.Range("D10").Value = "23 jun 2023"
$ws.Range("E10").Value = 1000
$ws.Range("F10").Value = 0

# --- Row 11: "Service Conversation" - only the date label changes ---
$ws.Range("D11").Value = "23 jun 2023"

# --- Row 12: "Marketing Conversation" - only the date label changes ---
$ws.Range("D12").Value = "23 jun 2023"

# --- Row 13: "Utility Conversation" - only the date label changes ---
$ws.Range("D13").Value = "23 jun 2023"

# --- Row 14: "Authentication Conversation" - only the date label changes ---
$ws.Range("D14").Value = "23 jun 2023"

# --- Row 15: Subtotal value changes from 1000 to 1703 ---
$ws.Range("F15").Value = 1703

# Rows 16-23 (ESTIMATED TOTAL rows + footnotes + trailing blank) shift up
# unchanged and already match the target content/styles.
